# Weekly update: insert a new data row for "Haba" (Feria Lagunitas de Puerto
# Montt) ahead of the existing history, pushing the prior rows 147-163 down
# to 148-164 (dimension grows from A1:R163 to A1:R164).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 147:163 down by inserting a fresh row at 147.
$ws.Rows(147).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A147").Value = 4
$ws.Range("B147").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C147").Value = 'Los Lagos'
$ws.Range("D147").Value = 45194
$ws.Range("E147").Value = 10
$ws.Range("F147").Value = 100112026
$ws.Range("G147").Value = 'Haba'
$ws.Range("H147").Value = 'Sin especificar'
$ws.Range("I147").Value = 'Primera'
$ws.Range("J147").Value = 40
$ws.Range("K147").Value = 18000
$ws.Range("L147").Value = 18000
$ws.Range("M147").Value = 18000
$ws.Range("N147").Value = '$/saco 25 kilos'
$ws.Range("O147").Value = 'Provincia de Limarí'
$ws.Range("P147").Value = 720
$ws.Range("Q147").Value = 25
$ws.Range("R147").Value = 'Hortaliza'
